$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 577.2025
$ws.Range("J17").Value = 578.43054
$ws.Range("L17").Value = 1735.29162
$ws.Range("N17").Value = -2071.29162

$ws.Range("H40").Value = 2220
$ws.Range("I40").Value = 2275
$ws.Range("J40").Value = 2176
$ws.Range("K40").Value = 2275
$ws.Range("L40").Value = 2176
$ws.Range("M40").Value = -2100
$ws.Range("N40").Value = -2526

$ws.Range("H51").Value = 2333.3333
$ws.Range("J51").Value = 2333.3333
$ws.Range("L51").Value = 2333.3333
$ws.Range("N51").Value = -3301.3333

$ws.Range("H64").Value = 3333.3333

$ws.Range("H67").Value = 3333.3333

$ws.Range("H113").Value = 2750
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 754
$ws.Range("N113").Value = -9508

$ws.Range("H138").Value = 2521.652
$ws.Range("J138").Value = 2589.7407
$ws.Range("L138").Value = 7769.222099999999
$ws.Range("N138").Value = -18049.2221

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1730.75
$ws.Range("I2").Value = 1835.3529
$ws.Range("K2").Value = 1835.3529
$ws.Range("M2").Value = -1722.3529

$ws.Range("H45").Value = 2194.818
$ws.Range("I45").Value = 2092.6924
$ws.Range("J45").Value = 2261.2
$ws.Range("K45").Value = 2092.6924
$ws.Range("L45").Value = 2261.2
$ws.Range("M45").Value = -1715.6924
$ws.Range("N45").Value = -3015.2

$ws.Range("H61").Value = 3807.8
$ws.Range("I61").Value = 3740.5
$ws.Range("J61").Value = 3852.6667
$ws.Range("K61").Value = 3740.5
$ws.Range("L61").Value = 3852.6667
$ws.Range("M61").Value = -3528.5
$ws.Range("N61").Value = -4276.6667

$ws.Range("H74").Value = 2708
$ws.Range("I74").Value = 1687.5
$ws.Range("K74").Value = 1687.5
$ws.Range("M74").Value = -813.5

$ws.Range("H77").Value = 2708
$ws.Range("I77").Value = 1687.5
$ws.Range("K77").Value = 8437.5
$ws.Range("M77").Value = -4069.5

$ws.Range("H102").Value = 2141.6667
$ws.Range("I102").Value = 2099.25
$ws.Range("J102").Value = 2990
$ws.Range("K102").Value = 2099.25
$ws.Range("L102").Value = 2990
$ws.Range("M102").Value = -477.25
$ws.Range("N102").Value = -6234

$ws.Range("H106").Value = 41111
$ws.Range("J106").Value = 41111
$ws.Range("L106").Value = 41111
$ws.Range("N106").Value = -43635

$ws.Range("H110").Value = 3000
$ws.Range("I110").Value = 2000
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 2000
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = 45
$ws.Range("N110").Value = -9090

$ws.Range("H116").Value = 1730.75
$ws.Range("I116").Value = 1835.3529
$ws.Range("K116").Value = 1835.3529
$ws.Range("M116").Value = 458.6470999999999

$ws.Range("H132").Value = 4984.55
$ws.Range("I132").Value = 5196.706
$ws.Range("J132").Value = 4827.7393
$ws.Range("K132").Value = 15590.118
$ws.Range("L132").Value = 14483.2179
$ws.Range("M132").Value = -13060.118
$ws.Range("N132").Value = -19543.2179

$ws.Range("H136").Value = 3807.8
$ws.Range("I136").Value = 3740.5
$ws.Range("J136").Value = 3852.6667
$ws.Range("K136").Value = 11221.5
$ws.Range("L136").Value = 11558.0001
$ws.Range("M136").Value = -8671.5
$ws.Range("N136").Value = -16658.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1730.75
$ws.Range("I3").Value = 1835.3529
$ws.Range("K3").Value = 1835.3529
$ws.Range("M3").Value = -1721.3529

$ws.Range("H51").Value = 72000
$ws.Range("J51").Value = 72000
$ws.Range("L51").Value = 72000
$ws.Range("N51").Value = -72982

$ws.Range("H101").Value = 100000
$ws.Range("J101").Value = 100000
$ws.Range("L101").Value = 100000
$ws.Range("N101").Value = -106490

$ws.Range("H105").Value = 20836044
$ws.Range("I105").Value = 31252742
$ws.Range("J105").Value = 2650
$ws.Range("K105").Value = 31252742
$ws.Range("L105").Value = 2650
$ws.Range("M105").Value = -31250995
$ws.Range("N105").Value = -6144

$ws.Range("H134").Value = 3308.1035
$ws.Range("I134").Value = 3081.9375
$ws.Range("J134").Value = 3586.4614
$ws.Range("K134").Value = 9245.8125
$ws.Range("L134").Value = 10759.3842
$ws.Range("M134").Value = -6710.8125
$ws.Range("N134").Value = -15829.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1800.8
$ws.Range("I58").Value = 1490.5
$ws.Range("J58").Value = 2155.4285
$ws.Range("K58").Value = 1490.5
$ws.Range("L58").Value = 2155.4285
$ws.Range("M58").Value = -1287.5
$ws.Range("N58").Value = -2561.4285

$ws.Range("H117").Value = 49400
$ws.Range("J117").Value = 49400
$ws.Range("L117").Value = 49400
$ws.Range("N117").Value = -58578

$ws.Range("H122").Value = 1416.7354
$ws.Range("I122").Value = 945.6
$ws.Range("K122").Value = 2836.8
$ws.Range("M122").Value = -386.8000000000002

$ws.Range("H136").Value = 1800.8
$ws.Range("I136").Value = 1490.5
$ws.Range("J136").Value = 2155.4285
$ws.Range("K136").Value = 4471.5
$ws.Range("L136").Value = 6466.2855
$ws.Range("M136").Value = -1921.5
$ws.Range("N136").Value = -11566.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1275
$ws.Range("J31").Value = 1275
$ws.Range("L31").Value = 3825
$ws.Range("N31").Value = -4401

$ws.Range("H33").Value = 11198.444
$ws.Range("I33").Value = 33364.332
$ws.Range("J33").Value = 115.5
$ws.Range("K33").Value = 200185.992
$ws.Range("L33").Value = 693
$ws.Range("M33").Value = -199902.992
$ws.Range("N33").Value = -1259

$ws.Range("H35").Value = 300
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5716.8823
$ws.Range("I70").Value = 5725.727
$ws.Range("J70").Value = 5700.6665
$ws.Range("K70").Value = 5725.727
$ws.Range("L70").Value = 5700.6665
$ws.Range("M70").Value = -5455.727
$ws.Range("N70").Value = -6240.6665

$ws.Range("H73").Value = 5716.8823
$ws.Range("I73").Value = 5725.727
$ws.Range("J73").Value = 5700.6665
$ws.Range("K73").Value = 5725.727
$ws.Range("L73").Value = 5700.6665
$ws.Range("M73").Value = -4789.727
$ws.Range("N73").Value = -7572.6665

$ws.Range("H102").Value = 1569.75
$ws.Range("I102").Value = 1416.25
$ws.Range("J102").Value = 1800
$ws.Range("K102").Value = 1416.25
$ws.Range("L102").Value = 1800
$ws.Range("M102").Value = 205.75
$ws.Range("N102").Value = -5044

$ws.Range("H113").Value = 1265.25
$ws.Range("J113").Value = 1397.1666
$ws.Range("L113").Value = 1397.1666
$ws.Range("N113").Value = -5737.1666

$ws.Range("H126").Value = 2196
$ws.Range("I126").Value = 2196
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6588
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -4118

$ws.Range("H132").Value = 3130.6667
$ws.Range("I132").Value = 3387.4285
$ws.Range("J132").Value = 2967.2727
$ws.Range("K132").Value = 10162.2855
$ws.Range("L132").Value = 8901.8181
$ws.Range("M132").Value = -7632.2855
$ws.Range("N132").Value = -13961.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 128749.875
$ws.Range("I40").Value = 204000.8
$ws.Range("J40").Value = 3331.6667
$ws.Range("K40").Value = 204000.8
$ws.Range("L40").Value = 3331.6667
$ws.Range("M40").Value = -203864.8
$ws.Range("N40").Value = -3603.6667

$ws.Range("H61").Value = 7180.4
$ws.Range("I61").Value = 6614.857
$ws.Range("J61").Value = 8500
$ws.Range("K61").Value = 6614.857
$ws.Range("L61").Value = 8500
$ws.Range("M61").Value = -6412.857
$ws.Range("N61").Value = -8904

$ws.Range("H82").Value = 23812698
$ws.Range("I82").Value = 41670036
$ws.Range("J82").Value = 2911.111
$ws.Range("K82").Value = 41670036
$ws.Range("L82").Value = 2911.111
$ws.Range("M82").Value = -41669675
$ws.Range("N82").Value = -3633.111

$ws.Range("H85").Value = 23812698
$ws.Range("I85").Value = 41670036
$ws.Range("J85").Value = 2911.111
$ws.Range("K85").Value = 41670036
$ws.Range("L85").Value = 2911.111
$ws.Range("M85").Value = -41668788
$ws.Range("N85").Value = -5407.111

$ws.Range("H113").Value = 7180.4
$ws.Range("I113").Value = 6614.857
$ws.Range("J113").Value = 8500
$ws.Range("K113").Value = 6614.857
$ws.Range("L113").Value = 8500
$ws.Range("M113").Value = -4444.857
$ws.Range("N113").Value = -12840

$ws.Range("H122").Value = 3142.2173
$ws.Range("I122").Value = 2981.2
$ws.Range("K122").Value = 8943.599999999999
$ws.Range("M122").Value = -6493.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4641.271
$ws.Range("I96").Value = 2848.4167
$ws.Range("K96").Value = 2848.4167
$ws.Range("M96").Value = -1475.4167

$ws.Range("H113").Value = 1389.9333
$ws.Range("I113").Value = 1701.909
$ws.Range("J113").Value = 532
$ws.Range("K113").Value = 5105.727000000001
$ws.Range("L113").Value = 1596
$ws.Range("M113").Value = -2935.727000000001
$ws.Range("N113").Value = -5936

$ws.Range("H122").Value = 2570.36
$ws.Range("I122").Value = 2022.4375
$ws.Range("J122").Value = 3544.4443
$ws.Range("K122").Value = 6067.3125
$ws.Range("L122").Value = 10633.3329
$ws.Range("M122").Value = -3617.3125
$ws.Range("N122").Value = -15533.3329
